# Update 16 May: add new uploads
# Adds a new row (A3) containing the shared string "sdfsf" below the
# existing "આપનું નામ" / "hkjhkhkjh" entries, which:
#   - appends "sdfsf" to xl/sharedStrings.xml (new shared string index 2)
#   - adds <row r="3">...</row> with cell A3 referencing that string
#   - extends the sheet dimension from A1:A2 to A1:A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "sdfsf"
